$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B28").Value = 39232
$ws.Range("B29").Value = 39262
$ws.Range("C29").Value = 668
$ws.Range("D29").Value = 5089.1099999999997
$ws.Range("B30").Value = 39284
$ws.Range("C30").Formula = "=B30-`$B`$6"
$ws.Range("D30").Value = 5142.46

$ws.Range("B7:D7").Copy()
$ws.Range("B30:D30").PasteSpecial(-4122)
